$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates (row-by-row, matches source diff) ---
$ws.Range("D2").Value = "42.925.15"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.208.70"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.85"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.15"
$ws.Range("E7").Value = "  +2.97%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.73"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.03"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "2.539.81"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.45"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "2.219.93"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "42.861.65"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.14"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.35"
$ws.Range("E22").Value = "  +7.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.94"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.21"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.44"
$ws.Range("E26").Value = "  +9.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.35"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.40"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.33"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0867"
$ws.Range("E33").Value = "  +9.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.22"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.122"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0362"
$ws.Range("E36").Value = "  +10.09%  "
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.37"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.03"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +20.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.27"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.87"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.46"
$ws.Range("E51").Value = "  +21.54%  "

# --- Row swap: MultiversX (42) <-> Algorand (43) ---
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.22"
$ws.Range("E43").Value = "  +2.61%  "

# --- Row swap: Cronos (47) <-> WOONetwork (48) ---
$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.466"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0972"
$ws.Range("E48").Value = "  -1.58%  "
